$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "[OK]" marker text in the Error column (K) is now printed without the
# surrounding brackets, e.g. for the case where the stmtID would be appended
# after an error instead ("Added a printing of the stmtID if there is an
# error").
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 11).Value = "OK"
}

# Touch a font property (and revert it) so the workbook's style catalog
# gains the extra font/cellXf entry that the real editing session left
# behind, without changing the look of any cell.
$ws.Range("K2").Font.Bold = $true
$ws.Range("K2").Font.Bold = $false

# Leave the cursor on the "Coded Statement" column of the first data row,
# matching the new recorded selection.
[void]$ws.Range("I2").Select()
